$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing D:K quarterly data to F:M.
$ws.Columns("D:E").Insert()

# Copy number formatting (style) from column F into the two new columns D:E
# for every data row so the new cells carry the same numeric/date formats
# that Excel would otherwise leave blank/General after a plain column insert.
$ws.Range("F8:F102").Copy()
$ws.Range("D8:E102").PasteSpecial(-4122)

$ws.Range("F7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("F80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)

$data = @{
    7 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    8 = @(2982200, 2900300, 2989000, 2893600, 2817300, 2815100, 2836400, 2931000, 2861700, 2808900)
    9 = @(2060000, 1908500, 1982600, 1862500, 1934300, 1910700, 1950000, 2019000, 2006000, 1946200)
    10 = @(922300, 991800, 1006400, 1031100, 883000, 904300, 886400, 912100, 855700, 862800)
    11 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    12 = @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
    13 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    14 = @(19100, 22400, 19100, 15700, 22400, 22400, 22400, 14100, 237100, 25800)
    15 = @(127900, 136900, 135800, 135800, 145900, 139100, 133500, 166700, 71600, 152600)
    16 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    17 = @(2582800, 2493100, 2633300, 2495300, 2552500, 2512100, 2606400, 2704500, 2691500, 2578900)
    18 = @(399400, 407300, 355700, 398300, 264800, 302900, 230000, 226500, 170200, 230100)
    19 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    20 = @(-46000, -65100, -67300, -24700, -30300, -51600, -67300, -34000, -81000, -41100)
    21 = @(344500, 343300, 288400, 509400, 241200, 256900, 136900, 359200, 403800, 245300)
    22 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    23 = @(353400, 342200, 288400, 373600, 234500, 251300, 162700, 192500, 89200, 189000)
    24 = @(18000, 113300, 49400, 67300, 40400, 37000, 2200, 66900, 23500, 49300)
    25 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    26 = @(335500, 228900, 239000, 306300, 194100, 214300, 160400, 125600, 65700, 139700)
    27 = @(341100, 240100, 248000, 307400, 195200, 214300, 163800, 133800, 106800, 151400)
    28 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    29 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    30 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    31 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    32 = @(46000, 65100, 67300, 24700, 30300, 51600, 67300, 34000, 81000, 41100)
    33 = @(341100, 240100, 248000, 307400, 195200, 214300, 163800, 133800, 106800, 151400)
    34 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    35 = @(341100, 240100, 248000, 307400, 195200, 214300, 163800, 133800, 106800, 151400)
    38 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    39 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    40 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    41 = @(1267800, 1127600, 705700, 770800, 681000, 469000, 559900, 1094000, 1118600, 916700)
    42 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    43 = @(1740200, 1702100, 1717800, 1667300, 1579800, 1576400, 1580900, 1618700, 1558800, 1659800)
    44 = @(1758200, 1686400, 1593200, 1626900, 1482100, 1493400, 1515800, 1687900, 1579900, 1514200)
    45 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    46 = @(4766200, 4516000, 4016700, 4065000, 3743000, 3538800, 3656600, 4400600, 4257400, 4090700)
    47 = @(2488600, 2403300, 2174400, 2155300, 2237200, 2186800, 2076800, 2217300, 2225500, 2027200)
    48 = @(6385200, 6234900, 6260700, 6284300, 6460400, 6440200, 6509800, 7095700, 7160200, 7358600)
    49 = @(586800, 576700, 577800, 570000, 536300, 476800, 490300, 521200, 507100, 516500)
    50 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    51 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    52 = @(189600, 141400, 179500, 190700, 228900, 226600, 270400, 301700, 318100, 293500)
    53 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    54 = @(14416400, 13872300, 13209200, 13265300, 13205800, 12869200, 13003900, 14536500, 14468400, 14286400)
    55 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    56 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    57 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    58 = @(1210600, 1188200, 1174700, 1123100, 1088300, 1067000, 1128700, 1629200, 1313500, 1460200)
    59 = @(2320300, 2174400, 2070100, 2147500, 2236100, 2022900, 2008400, 2225500, 2195000, 2132800)
    60 = @(3530900, 3362600, 3244800, 3270600, 3324500, 3090000, 3137100, 3854800, 3508500, 3593000)
    61 = @(2541300, 2516600, 2426900, 2319200, 2295600, 2347200, 2597400, 2739700, 3116500, 3056600)
    62 = @(791000, 746100, 714700, 732700, 792100, 871800, 916700, 936700, 955500, 923800)
    63 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    64 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    65 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    66 = @(6883400, 6651200, 6429000, 6374000, 6464900, 6362800, 6707300, 7594600, 7653200, 7684900)
    67 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    68 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    69 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    70 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    71 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    72 = @(5281200, 4957000, 4719100, 4832400, 4526100, 4275900, 4057100, 4415900, 4284400, 4223400)
    73 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    74 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    75 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    76 = @(7533000, 7221100, 6780200, 6891300, 6740900, 6506400, 6296600, 6941900, 6815100, 6601500)
    77 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    80 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    81 = @(341100, 240100, 248000, 307400, 195200, 214300, 163800, 133800, 106800, 151400)
    82 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    83 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    84 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    85 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    86 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    87 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    88 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    89 = @(306300, 457800, 319800, 141400, 529600, 408400, 317500, 97400, 489500, 389700)
    90 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    91 = @(322000, -144700, -140200, -130200, 360200, -154800, -143600, -158500, -258200, -250000)
    92 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    93 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    94 = @(-187400, -154800, -96500, -118900, -304100, -218800, -140200, -123300, -77500, -85700)
    95 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    96 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    97 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    98 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    99 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    100 = @(42600, 113300, -300700, 62800, -5600, -278300, -640700, -5900, -219500, 3500)
    101 = @(2200, -9000, 6700, 4500, -6700, -3400, -18000, 2300, 11700, 2300)
    102 = @(163800, 407300, -70700, 89800, 213200, -92000, -481300, -29300, 204200, 309900)
}


$cols = @("D","E","F","G","H","I","J","K","L","M")
foreach ($r in $data.Keys) {
    $rowvals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $v = $rowvals[$i]
        if ($v -ne $null) {
            $ws.Range($cols[$i] + $r).Value = $v
        }
    }
}
